$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 0
$ws.Range("I12").Value = 0
$ws.Range("K12").Value = 0
$ws.Range("M12").ClearContents()

$ws.Range("H33").Value = 280.64285
$ws.Range("I33").Value = 253.9
$ws.Range("K33").Value = 253.9
$ws.Range("M33").Value = -24.90000000000001

$ws.Range("H43").Value = 14999.75
$ws.Range("I43").Value = 5000
$ws.Range("K43").Value = 5000
$ws.Range("M43").Value = -4931

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 12502469
$ws.Range("J2").Value = 3442.75
$ws.Range("L2").Value = 3442.75
$ws.Range("N2").Value = -3668.75

$ws.Range("H61").Value = 9299.333000000001
$ws.Range("I61").Value = 9299.333000000001
$ws.Range("K61").Value = 9299.333000000001
$ws.Range("M61").Value = -9087.333000000001

$ws.Range("H80").Value = 99993.5
$ws.Range("J80").Value = 99993.5
$ws.Range("L80").Value = 99993.5
$ws.Range("N80").Value = -101989.5

$ws.Range("H83").Value = 99993.5
$ws.Range("J83").Value = 99993.5
$ws.Range("L83").Value = 299980.5
$ws.Range("N83").Value = -309964.5

$ws.Range("H88").Value = 2324.3333
$ws.Range("I88").Value = 1987
$ws.Range("J88").Value = 2493
$ws.Range("K88").Value = 1987
$ws.Range("L88").Value = 2493
$ws.Range("M88").Value = -1581
$ws.Range("N88").Value = -3305

$ws.Range("H91").Value = 2324.3333
$ws.Range("I91").Value = 1987
$ws.Range("J91").Value = 2493
$ws.Range("K91").Value = 1987
$ws.Range("L91").Value = 2493
$ws.Range("M91").Value = -583
$ws.Range("N91").Value = -5301

$ws.Range("H102").Value = 1733.6522
$ws.Range("I102").Value = 1251.7333
$ws.Range("K102").Value = 1251.7333
$ws.Range("M102").Value = 370.2666999999999

$ws.Range("H110").Value = 2791
$ws.Range("I110").Value = 2428.3
$ws.Range("K110").Value = 2428.3
$ws.Range("M110").Value = -383.3000000000002

$ws.Range("H116").Value = 12502469
$ws.Range("J116").Value = 3442.75
$ws.Range("L116").Value = 3442.75
$ws.Range("N116").Value = -8030.75

$ws.Range("H124").Value = 75000
$ws.Range("J124").Value = 75000
$ws.Range("L124").Value = 75000
$ws.Range("N124").Value = -84820

$ws.Range("H125").Value = 70000
$ws.Range("J125").Value = 70000
$ws.Range("L125").Value = 70000
$ws.Range("N125").Value = -79840

$ws.Range("H132").Value = 3590
$ws.Range("I132").Value = 2327.5
$ws.Range("K132").Value = 6982.5
$ws.Range("M132").Value = -4452.5

$ws.Range("H136").Value = 9299.333000000001
$ws.Range("I136").Value = 9299.333000000001
$ws.Range("K136").Value = 27897.999
$ws.Range("M136").Value = -25347.999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 12502469
$ws.Range("J3").Value = 3442.75
$ws.Range("L3").Value = 3442.75
$ws.Range("N3").Value = -3670.75

$ws.Range("H99").Value = 1704.8667
$ws.Range("I99").Value = 1911.9
$ws.Range("K99").Value = 1911.9
$ws.Range("M99").Value = -413.9000000000001

$ws.Range("H105").Value = 3128.0967
$ws.Range("I105").Value = 2810.3333
$ws.Range("K105").Value = 2810.3333
$ws.Range("M105").Value = -1063.3333

$ws.Range("H107").Value = 5999.5
$ws.Range("I107").Value = 0
$ws.Range("J107").Value = 5999.5
$ws.Range("K107").Value = 0
$ws.Range("L107").Value = 5999.5
$ws.Range("N107").Value = -9839.5
$ws.Range("M107").ClearContents()

$ws.Range("H134").Value = 1438
$ws.Range("I134").Value = 747.6667
$ws.Range("K134").Value = 2243.0001
$ws.Range("M134").Value = 291.9998999999998

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4709.75
$ws.Range("I31").Value = 3711.8572
$ws.Range("J31").Value = 5485.8887
$ws.Range("K31").Value = 3711.8572
$ws.Range("L31").Value = 5485.8887
$ws.Range("M31").Value = -3416.8572
$ws.Range("N31").Value = -6075.8887

$ws.Range("H34").Value = 4709.75
$ws.Range("I34").Value = 3711.8572
$ws.Range("J34").Value = 5485.8887
$ws.Range("K34").Value = 3711.8572
$ws.Range("L34").Value = 5485.8887
$ws.Range("M34").Value = -3509.8572
$ws.Range("N34").Value = -5889.8887

$ws.Range("H62").Value = 197254.5
$ws.Range("I62").Value = 129673
$ws.Range("K62").Value = 129673
$ws.Range("M62").Value = -129049

$ws.Range("H65").Value = 197254.5
$ws.Range("I65").Value = 129673
$ws.Range("K65").Value = 648365
$ws.Range("M65").Value = -645245

$ws.Range("H141").Value = 93333.336
$ws.Range("J141").Value = 93333.336
$ws.Range("L141").Value = 93333.336
$ws.Range("N141").Value = -103693.336

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H6").Value = 88.57143000000001
$ws.Range("I6").Value = 70
$ws.Range("K6").Value = 210
$ws.Range("M6").Value = -97

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H33").Value = 16750
$ws.Range("J33").Value = 16750
$ws.Range("L33").Value = 16750
$ws.Range("N33").Value = -17254

$ws.Range("H70").Value = 4873.75
$ws.Range("I70").Value = 4747.5
$ws.Range("K70").Value = 4747.5
$ws.Range("M70").Value = -4477.5

$ws.Range("H73").Value = 4873.75
$ws.Range("I73").Value = 4747.5
$ws.Range("K73").Value = 4747.5
$ws.Range("M73").Value = -3811.5

$ws.Range("H92").Value = 13873.333
$ws.Range("J92").Value = 13873.333
$ws.Range("L92").Value = 13873.333
$ws.Range("N92").Value = -17617.333

$ws.Range("H97").Value = 2889.2222
$ws.Range("J97").Value = 3300.7
$ws.Range("L97").Value = 3300.7
$ws.Range("N97").Value = -4292.7

$ws.Range("H102").Value = 3148
$ws.Range("I102").Value = 1949.5
$ws.Range("K102").Value = 1949.5
$ws.Range("M102").Value = -327.5

$ws.Range("H113").Value = 31287136
$ws.Range("I113").Value = 62519450
$ws.Range("J113").Value = 54822.25
$ws.Range("K113").Value = 62519450
$ws.Range("L113").Value = 54822.25
$ws.Range("M113").Value = -62517280
$ws.Range("N113").Value = -59162.25

$ws.Range("H123").Value = 34999.832
$ws.Range("I123").Value = 0
$ws.Range("J123").Value = 34999.832
$ws.Range("K123").Value = 0
$ws.Range("L123").Value = 34999.832
$ws.Range("N123").Value = -39899.832
$ws.Range("M123").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 5999.8
$ws.Range("I132").Value = 6499.6665
$ws.Range("K132").Value = 19498.9995
$ws.Range("M132").Value = -16968.9995

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H29").Value = 0
$ws.Range("I29").Value = 0
$ws.Range("K29").Value = 0
$ws.Range("M29").ClearContents()

$ws.Range("H132").Value = 2443.889
$ws.Range("I132").Value = 2249.375
$ws.Range("K132").Value = 6748.125
$ws.Range("M132").Value = -4218.125

$ws.Range("H136").Value = 902.8333
$ws.Range("I136").Value = 720.64703
$ws.Range("K136").Value = 2161.94109
$ws.Range("M136").Value = 388.0589100000002
